# Fill in zero placeholders for the balance-sheet data bucketing columns
# on the three statement sheets (cbs_7, cpl_6, ccf_9). These cells were
# previously blank and are now explicitly set to 0.

$wb = $excel.ActiveWorkbook

# --- Sheet: cbs_7 (balance sheet) ---
$ws1 = $wb.Worksheets.Item("cbs_7")

$ws1.Range("C2").Value = 0
$ws1.Range("D2").Value = 0

$ws1.Range("C3").Value = 0
$ws1.Range("D3").Value = 0

$ws1.Range("C7").Value = 0

$ws1.Range("C11").Value = 0
$ws1.Range("D11").Value = 0

$ws1.Range("C15").Value = 0

$ws1.Range("C17").Value = 0

$ws1.Range("C20").Value = 0
$ws1.Range("D20").Value = 0

$ws1.Range("C21").Value = 0
$ws1.Range("D21").Value = 0

$ws1.Range("C23").Value = 0

$ws1.Range("D25").Value = 0

$ws1.Range("C28").Value = 0
$ws1.Range("D28").Value = 0

$ws1.Range("C36").Value = 0
$ws1.Range("D36").Value = 0

# --- Sheet: cpl_6 (profit & loss) ---
$ws2 = $wb.Worksheets.Item("cpl_6")

$ws2.Range("D12").Value = 0

# --- Sheet: ccf_9 (cash flow) ---
$ws3 = $wb.Worksheets.Item("ccf_9")

$ws3.Range("C2").Value = 0
$ws3.Range("D2").Value = 0

$ws3.Range("C10").Value = 0
$ws3.Range("D10").Value = 0

$ws3.Range("C13").Value = 0
$ws3.Range("D13").Value = 0

$ws3.Range("C15").Value = 0
